$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row: "_old" -> "_FV2310" (columns A-J) and "_new" -> "_FV2404" (columns L-U)
$oldCols = @("A","B","C","D","E","F","G","H","I","J")
foreach ($c in $oldCols) {
    $cell = $ws.Range($c + "1")
    $cell.Value = ($cell.Value2 -replace "_old$", "_FV2310")
}

$newCols = @("L","M","N","O","P","Q","R","S","T","U")
foreach ($c in $newCols) {
    $cell = $ws.Range($c + "1")
    $cell.Value = ($cell.Value2 -replace "_new$", "_FV2404")
}

# Freeze the header row (row 1)
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# Turn the data range into a table
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U56"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
